$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 13 data: GFG | Print sum of all subset | Java | 2023-02-26
$ws.Range("A13").Value = "GFG"
$ws.Range("B13").Value = "Print sum of all subset"
$ws.Range("C13").Value = "Java"

# D13 needs the date-formatted style like D11/D12 (numFmtId 15)
$ws.Range("D12").Copy()
$ws.Range("D13").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("D13").Value = 44983

# Update selection to D14 as per diff
$ws.Range("D14").Select()
